# Protocollo CS Chat.xlsx - applies the authored edit:
#  - updates the "msgs§..." protocol string (C4) to add <hour>/<minute> fields
#    and highlights it with a yellow fill
#  - adds a new empty, underlined, selected cell at G3
#  - tidies row heights / page setup the way Excel does when the file is
#    re-saved

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C4: new protocol string for "msgs", now carrying <hour> and <minute> ---
$ws.Range("C4").Value2 = "msgs§<<nickname>&<type>&<msg>&<hour>&<minute>>§<<nickname>&<type>&<msg>&<hour>&<minute>>§<…>"
$ws.Range("C4").Interior.Color = 65535
$ws.Range("C4").HorizontalAlignment = -4108
$ws.Range("C4").WrapText = $true

# --- G3: new empty, underlined cell, left as the active selection ---
$ws.Range("G3").Font.Underline = 2
$ws.Range("G3").Select()

# --- row heights (re-flowed by Excel on save) ---
$ws.Rows.Item(1).RowHeight = 18.75
$ws.Rows.Item(2).RowHeight = 126
$ws.Rows.Item(3).RowHeight = 47.25
$ws.Rows.Item(4).RowHeight = 157.5
$ws.Rows.Item(5).RowHeight = 110.25
$ws.Rows.Item(6).RowHeight = 63
$ws.Rows.Item(7).RowHeight = 47.25
$ws.Rows.Item(8).RowHeight = 47.25
$ws.Rows.Item(9).RowHeight = 63

# --- page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "done"
